# Modules.xlsx update: remove "ElementName3"/col F, and fill in Name/ChefModule/
# ElementName1/ElementName2 details for each module row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to have an extra (always empty) column F with header
# "ElementName3". Drop that whole column first - this also removes the
# now-unused "ElementName3" shared string and shifts everything else left.
$ws.Columns("F").Delete() | Out-Null

# Fill in the per-module details that used to be blank.
$ws.Range("B2").Value = "pede. Suspendisse dui."
$ws.Range("C2").Value = "EL Haddad"
$ws.Range("D2").Value = "Nullam feugiat placerat"
$ws.Range("E2").Value = "varius et, euismod"

$ws.Range("B3").Value = "a nunc. In"
$ws.Range("C3").Value = "Badir"
$ws.Range("D3").Value = "sodales nisi magna"
$ws.Range("E3").Value = "elementum sem, vitae"

$ws.Range("B4").Value = "amet metus. Aliquam"
$ws.Range("C4").Value = "Ezzine"
$ws.Range("D4").Value = "Cras vulputate velit"
$ws.Range("E4").Value = "scelerisque neque sed"

$ws.Range("B5").Value = "quam vel sapien"
$ws.Range("C5").Value = "El Alami Hassoun"
$ws.Range("D5").Value = "Nunc mauris elit,"
$ws.Range("E5").Value = "libero et tristique"

$ws.Range("B6").Value = "feugiat nec, diam."
$ws.Range("C6").Value = "Lazaar"
$ws.Range("D6").Value = "pellentesque. Sed dictum."
$ws.Range("E6").Value = "ridiculus mus. Proin"

$ws.Range("B7").Value = "nonummy. Fusce fermentum"
$ws.Range("C7").Value = "El Haddad"
$ws.Range("D7").Value = "neque pellentesque massa"
$ws.Range("E7").Value = "Mauris eu turpis."

$ws.Range("B8").Value = "a, arcu. Sed"
$ws.Range("C8").Value = "EL Haddad"
$ws.Range("D8").Value = "sit amet risus."
$ws.Range("E8").Value = "Nulla facilisi. Sed"

$ws.Range("B9").Value = "Suspendisse eleifend. Cras"
$ws.Range("C9").Value = "El Alami Hassoun"
$ws.Range("D9").Value = "velit dui, semper"
$ws.Range("E9").Value = "ligula elit, pretium"

# B2 picked up an explicit (non-theme) black font colour along the way.
$ws.Range("B2").Font.Color = 0

# Re-fit the columns to their new, wider content.
$ws.Columns("A").ColumnWidth = 8.736979166666666
$ws.Columns("B").ColumnWidth = 26.166666666666668
$ws.Columns("C").ColumnWidth = 15.307291666666666
$ws.Columns("D").ColumnWidth = 24.307291666666668
$ws.Columns("E").ColumnWidth = 20.451822916666668

# Leave the selection where the user ended up after entering the data.
$ws.Range("F11").Select() | Out-Null
